$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51
$ws.Range("A45").Copy()
$ws.Range("A51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A51").Value = 43881
$ws.Range("B51").Value = "5pm - 7.50 pm"
$ws.Range("C51").Value = "N/A"
$ws.Range("D51").Value = "Was expecting mid-term results, looking forward to know more about the system architecture"
$ws.Range("E51").Value = "Learned three new Key Expert practices.`nLearned about the architecture of the system and the process of understanding the architecture from source code."
$ws.Range("F51").Value = "If there's no proper documented architecture, the first step in understanding the architecture from the source code is by looking at the folder/package  structure and trying to grouping various related classes together. We can start with the UML diagram and slowly try to abstract up. `nPull requests can also be a useful tool as lot of design decisions could be present in pull requests which will help us to understand the rationale of the developer and why certain things are written in a certain way. "
$ws.Range("G51").Value = "This week's homework is pretty vast and have to start working on it soon."
$ws.Rows(51).RowHeight = 249.6

# Row 53
$ws.Range("A45").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A53").Value = 43884
$ws.Range("B53").Value = "1pm - 6pm"
$ws.Range("C53").Value = "Vaishakhi, Aman"
$ws.Range("D53").Value = "Understand the architecture of the system and document the same"
$ws.Range("E53").Value = "Understood both the as-described and as-implemented architecture of the system."
$ws.Range("F53").Value = "Realized the importance of having a documented version of architecture as it helps developers in the future `nThere are only few variations in the as-described and as-implemented architectures of h2, which means the h2 community has done a really good job at code reviews and maintaining the standards`n"
$ws.Range("G53").Value = "Proud to complete the most challenging part of this week's assignment"
$ws.Rows(53).RowHeight = 156

# Row 55
$ws.Range("A45").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A55").Value = 43885
$ws.Range("B55").Value = "9pm - 12am"
$ws.Range("C55").Value = "Vaishakhi, Aman"
$ws.Range("D55").Value = "Finish the remaining sections of the homework like pull requests, issues, state of the system etc. "
$ws.Range("E55").Value = "We divided the remaining sections of the homework and worked on it individually. Later on we collated all our findings in the report."
$ws.Range("F55").Value = "Looking for the social context was not difficult as h2 maintains good documentation in their website`nUnderstood the importance of having proper comments and explanation in the pull requests as it helped in understanding the decisions made by the contributors and the rationale behind each change"
$ws.Range("G55").Value = "Happy to complete the homework early"
$ws.Rows(55).RowHeight = 156

# Final view state: selection + zoom (matches author's saved view)
$excel.ActiveWindow.Zoom = 100
$ws.Range("F55").Select()
